# "asta reportes FINAL CLASE LUNES 28"
#
# The sheet originally holds a single shared-string value ("Hello World !")
# in A1. The edit renames that text to "LISTA DE PRECIOS !" and relocates it
# to C3 (the sheet's used range grows from A1:A1 to A1:C3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the updated text into its new home, C3 ...
$ws.Range("C3").Value = "LISTA DE PRECIOS !"

# ... and remove it from its old location, A1, so the cell is vacated.
$ws.Range("A1").ClearContents()
